# Fix: histogram-feeding columns (R/S/Y) were stale after the underlying
# dataset was reloaded; refresh the computed noise/detection-limit values
# for every data row (2-54) so the histogram reflects the reloaded data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (col index 18)
$valsR = @(0.1150564579126169,0.1329060090689602,0.1504895836599083,0.1495409269050694,0.1493590408051323,0.1414091646652366,0.1570249706070675,0.1378712521730261,0.1611859465070608,0.1699465612478473,0.1488961910983873,0.1729452700183139,0.1584492338038465,0.148128985927394,0.1776781307107918,0.1813154685736886,0.2216613326618502,0.2719872017944054,0.2362831439242929,0.1791845798329176,0.2899835967925228,0.1885154694327391,0.1998930229818494,0.2620628862971832,0.3111449374357326,0.2171764998636597,0.2481452745227085,0.2323032571743541,0.2336197477542778,0.2909845109322784,0.191236726899488,0.2147797867693591,0.2181191187883231,0.2486386107575973,0.2587354444252838,0.2210871590460066,0.2945624602804249,0.2282355956590671,0.3084664161991244,0.2817599513695137,0.2435329039394293,0.2164001039216811,0.2326939981581573,0.2820988899053841,0.2646564500047142,0.2720394608535724,0.2790928755616518,0.3074102121612757,0.3258284558540636,0.3616316432504626,0.3833230658380932,0.3683295506480807,0.3954497449951028)
for ($i = 0; $i -lt $valsR.Count; $i++) {
    $ws.Cells.Item(2 + $i, 18).Value = $valsR[$i]
}

# Column S (col index 19)
$valsS = @(0.7053430567844436,1.056417420870683,1.321067744810924,1.30703969975429,1.296420910492415,1.25939010822333,1.451793310954781,1.116966191759442,1.625515643193604,1.691119465268943,1.509974615348027,1.716543968242473,1.707838119149622,1.482214388592182,1.923568837070754,1.986890142742924,2.927404316009116,3.926925193376464,3.210037146657739,2.216357050357354,4.798587647809612,2.425999622904041,2.664662922801566,3.971009233802363,5.769462258847229,2.937274559609091,3.770312451142793,3.455316057831864,3.348055225913699,5.10700184009962,2.62211514883971,2.984498563121746,2.80885484021969,4.183112013683599,4.124248537728118,3.177148601375486,5.408479209716671,3.603824368200521,5.948115936703149,5.618113323764024,3.991353303058379,3.379204450581494,3.732921590247159,5.518356426012242,4.850626864704772,5.119217003608957,4.885914841066255,5.417884667322776,6.48109688732279,7.792020021145571,7.382133249675379,7.108863265882341,8.245614018096386)
for ($i = 0; $i -lt $valsS.Count; $i++) {
    $ws.Cells.Item(2 + $i, 19).Value = $valsS[$i]
}

# Column Y (col index 25)
$valsY = @(0.6059621814264952,0.699853783406671,0.7923265053645273,0.7873395523992832,0.7863789909894255,0.7445379750915352,0.826662952928614,0.7259687949280884,0.8484567408131748,0.8945301795413265,0.7838260976791516,0.9102999542921336,0.8339999521597034,0.7798117873511712,0.9350469101940364,0.9541275674381292,1.165135789197005,1.426673775128126,1.241373711142537,0.9427146303179958,1.516352172452704,0.9915794105479429,1.051095844497651,1.374409012390512,1.617799008302385,1.141495814244838,1.302084095200522,1.219858141800944,1.227053950332506,1.519275245929688,1.005647718073182,1.128844057565363,1.146704189146937,1.303106602682988,1.356270068699268,1.161566838121022,1.535180056971724,1.198095288795323,1.601433130223408,1.466367425926581,1.277546882433404,1.136704843935265,1.221857200049381,1.470521387698876,1.383787419862009,1.421226891655597,1.462069519986821,1.610545511730276,1.69693257191342,1.871107077411237,2.003039574922152,1.926677340081249,2.060467673532931)
for ($i = 0; $i -lt $valsY.Count; $i++) {
    $ws.Cells.Item(2 + $i, 25).Value = $valsY[$i]
}
